$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "62.776.36"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.12%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.054.63"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.27%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "538.44"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "134.44"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.66%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.053.35"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.50%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.492"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.55%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.153"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.28%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.19"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.34%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.450"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.07%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000222"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.20%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "34.09"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.40%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.533.89"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("E16").Value = "  +1.78%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "62.628.37"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.40%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.042.36"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.07%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.60"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.03%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "466.91"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.46%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.34"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.47%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.688"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.93%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "6.97"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.88%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "78.23"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.06"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("E26").Value = "  +0.03%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.69"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.65%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.80"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.29%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.47%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "25.92"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.95%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.15"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +4.87%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.87"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.21%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "57.64"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.06%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.46"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +5.88%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.28"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.10%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.92"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.75%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "461.78"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.56%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.204.97"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.36%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0391"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.16%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0793"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.53%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.117"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.32%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "8.09"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.51%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.53"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("E44").Value = "  +0.16%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.249"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.15%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "25.22"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.91%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "121.99"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.44%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.109"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.41%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.98"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.82%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0₃0516"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.29%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.25"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +5.67%  "
